$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E24").Value2 = "Variable  Importance using vip()"

$ws.Range("E25").Value2 = "Monthly: "
$ws.Range("H25").Value2 = "All Temporal Scale:"

$ws.Range("E26").Value2 = "              variable    Overall"
$ws.Range("H26").Value2 = "                                       variable    Overall"
$ws.Range("M26").Value2 = "Same"

$ws.Range("E27").Value2 = "1       srad_sum_April 6.79172717"
$ws.Range("H27").Value2 = "1                              gdd_sum_November  4.4389873"
$ws.Range("M27").Value2 = "tamp_mean_July"

$ws.Range("E28").Value2 = "2     tmax_mean_August 5.51094814"
$ws.Range("H28").Value2 = "2                                tamp_mean_July  4.2799157"
$ws.Range("M28").Value2 = "gdd_sum_November"

$ws.Range("E29").Value2 = "3       tamp_mean_July 5.37070104"
$ws.Range("H29").Value2 = "3                            gdd_mean_June_July  4.1157730"

$ws.Range("E30").Value2 = "4  dayl_mean_September 5.24681141"
$ws.Range("H30").Value2 = "4                        gdd_mean_growingseason  3.9875056"

$ws.Range("E31").Value2 = "5     gdd_sum_November 4.87230151"
$ws.Range("H31").Value2 = "5                            srad_sum_September  3.8805598"

$ws.Range("E24").Select()
